$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 10: "Max PE Silos Setting: " (user-defined) -------------------
# A10 is a label styled like the existing "Total Run Time" label (E5):
# bold-ish heading font at 18pt, thin border, light accent fill.
$a10 = $ws.Range("A10")
$a10.Value = "Max PE Silos Setting: "
$a10.Style = "Normal 4"

$a10Borders = $a10.Borders
$a10Borders.LineStyle = 1
$a10Borders.ColorIndex = 1

$a10.Interior.Pattern = 1
$a10.Interior.PatternColorIndex = -4105
$a10.Interior.ThemeColor = 8

$a10.Font.Size = 18

# B10 is the (empty) user input cell next to it, styled like the other
# input cells on the sheet (e.g. B5/B7/B8): 18pt font with a thin border.
$b10 = $ws.Range("B10")
$b10Borders = $b10.Borders
$b10Borders.LineStyle = 1
$b10Borders.ColorIndex = 1
$b10.Font.Size = 18

$ws.Rows.Item(10).RowHeight = 24

# The workbook was last left with this cell selected.
$ws.Range("B13").Select()
